# Switch example group numbers
# Replace the example reviewer names (Veselin, Rawda, Hannah, Mirit, Bogdana, Martin)
# back to the generic template names (Alice, Bob, Claire, David, Elaine), and remove
# the extra 6th template row (row 17) that is no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Restore the generic example names in column B (rows 12-16), clearing any
# custom font/color formatting those cells had picked up.
$ws.Range("B12").Value = "Alice"
$ws.Range("B13").Value = "Bob"
$ws.Range("B14").Value = "Claire"
$ws.Range("B15").Value = "David"
$ws.Range("B16").Value = "Elaine"

$ws.Range("B12:B16").Style = "Normal"

# Remove the sixth example row entirely (was B17/C17/E17 = Martin/3/2)
$ws.Rows.Item(17).Delete()

# Update the active selection to match the saved view state
$ws.Range("E14").Select()
